# Questionnaire edit: add a second "OG GEQ" table (columns G:L) mirroring
# the first table (columns A:E), plus GEQ code annotations in column L,
# and a new section label in A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New section label (set first so it lands at shared-string index 35,
#    matching the target workbook's shared string ordering).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 'OG GEQ'

# ---------------------------------------------------------------------
# 2. Populate the brand-new strings in the exact order they were first
#    authored, so the resulting shared-string table indices (36-44) line
#    up with the target workbook.
# ---------------------------------------------------------------------
$ws.Range("H15").Value = 'I was fully occupied with the game'
$ws.Range("L15").Value = 'GEQ33-5'
$ws.Range("L19").Value = 'GEQ33-28'
$ws.Range("H19").Value = 'I was deeply concentrated in the game'
$ws.Range("H11").Value = 'I was interested in the game''s story'
$ws.Range("L11").Value = 'GEQ33-3'
$ws.Range("H13").Value = 'I felt that I could explore things'
$ws.Range("H23").Value = 'I enjoyed it'
$ws.Range("L23").Value = 'GEQ33-20'

# ---------------------------------------------------------------------
# 3. Fill in the rest of the new G:K block, reusing existing question
#    text (and therefore existing shared-string entries) where the
#    question is identical to its A:E counterpart.
# ---------------------------------------------------------------------
$ws.Range("G9").Value = "Number"
$ws.Range("H9").Value = "Question"
$ws.Range("I9").Value = "No"
$ws.Range("J9").Value = "Maybe"
$ws.Range("K9").Value = "Yes"

$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 'I lose track of time'
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 'I feel different'
$ws.Range("G13").Value = 4
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 'The game feels real'
$ws.Range("G15").Value = 6
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 'I get wound up'
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 'Time seems to kind of stand still or stop'
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 'I feel spaced out'
$ws.Range("G19").Value = 10
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 'I can’t tell that I’m getting tired'
$ws.Range("G21").Value = 12
$ws.Range("H21").Value = 'Playing seems automatic'
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 'My thoughts go fast'
$ws.Range("G23").Value = 14
$ws.Range("G24").Value = 15
$ws.Range("H24").Value = 'I play without thinking about how to play'
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 'Playing makes me feel calm'
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 'I play longer than I meant to'
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = 'I really get into the game'
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = 'I feel like I just can’t stop playing'

# ---------------------------------------------------------------------
# 4. Apply styling to the new/changed cells by copying cell formats from
#    already-styled reference cells. This reproduces the same style
#    index (s="1" / s="2") rather than allocating brand-new ones.
# ---------------------------------------------------------------------
function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}


CopyFormat "A12" "G11"
CopyFormat "B12" "H11"
CopyFormat "C11" "I11"
CopyFormat "D11" "J11"
CopyFormat "E11" "K11"

CopyFormat "A12" "G12"
CopyFormat "B12" "H12"
CopyFormat "C12" "I12"
CopyFormat "D12" "J12"
CopyFormat "E12" "K12"

CopyFormat "A12" "G13"
CopyFormat "B12" "H13"
CopyFormat "C11" "I13"
CopyFormat "D11" "J13"
CopyFormat "E11" "K13"

CopyFormat "A12" "G14"
CopyFormat "B12" "H14"

CopyFormat "A12" "G15"
CopyFormat "B12" "H15"
CopyFormat "C11" "I15"
CopyFormat "D11" "J15"
CopyFormat "E11" "K15"

CopyFormat "A12" "G16"
CopyFormat "B12" "H16"

CopyFormat "A12" "G17"
CopyFormat "B12" "H17"

CopyFormat "A12" "G18"
CopyFormat "B12" "H18"

CopyFormat "A12" "G19"
CopyFormat "B12" "H19"
CopyFormat "C11" "I19"
CopyFormat "D11" "J19"
CopyFormat "E11" "K19"

# Row 20: A:E style changes from 1 to 2
CopyFormat "A12" "A20"
CopyFormat "B12" "B20"
CopyFormat "C12" "C20"
CopyFormat "D12" "D20"
CopyFormat "E12" "E20"
CopyFormat "A12" "G20"
CopyFormat "B12" "H20"
CopyFormat "C12" "I20"
CopyFormat "D12" "J20"
CopyFormat "E12" "K20"



# Row 23: A:E style changes from None to 1
CopyFormat "A11" "A23"
CopyFormat "B11" "B23"



# Row 26: A:E style changes from 1 to 2
CopyFormat "A12" "A26"
CopyFormat "B12" "B26"
CopyFormat "C12" "C26"
CopyFormat "D12" "D26"
CopyFormat "E12" "E26"
CopyFormat "A12" "G26"
CopyFormat "B12" "H26"
CopyFormat "C12" "I26"
CopyFormat "D12" "J26"
CopyFormat "E12" "K26"


# Row 28: A:E style changes from 1 to 2
CopyFormat "A12" "A28"
CopyFormat "B12" "B28"
CopyFormat "C12" "C28"
CopyFormat "D12" "D28"
CopyFormat "E12" "E28"
CopyFormat "A12" "G28"
CopyFormat "B12" "H28"
CopyFormat "C12" "I28"
CopyFormat "D12" "J28"
CopyFormat "E12" "K28"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Column widths for the new G/H columns (mirroring A/B).
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# ---------------------------------------------------------------------
# 6. Create the second table ("Tabela32") over G9:K28.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("G9:K28"), 0, 1)
$lo.Name = "Tabela32"
$lo.TableStyle = "TableStyleLight1"

# ---------------------------------------------------------------------
# 7. Final selection, matching the saved state in the target workbook.
# ---------------------------------------------------------------------
$ws.Range("F3").Select()
